# Apply the cryptos.xlsx price/volume/coin-swap update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.266.55"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "1.647.15"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.05%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "216.97"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("E9").Value = "  +0.47%  "
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.91"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.38%  "
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0793"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.875.95"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.29"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "1.645.46"
$ws.Range("E14").Value = "  +0.48%  "
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.547"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("E16").Value = "  +0.45%  "
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "63.38"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "26.266.77"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("E20").Value = "  -0.73%  "
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "195.72"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "143.48"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.78"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("E35").Value = "  +1.34%  "
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.913"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").Value = "1.138.35"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("E41").Value = "  -0.09%  "
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "100.54"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.18%  "
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.53"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "1.784.81"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  -0.77%  "
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "57.17"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +3.28%  "
$ws.Range("E48").Value = "  +3.05%  "
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.418"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.70"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +3.35%  "
